# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.47 = 17558.63 pesos`n✅ 17558.63 pesos = 4.44 = 942.8 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 223.85
$wsTasas.Range("O10").Value = 3930.5
$wsTasas.Range("N12").Value = 3957.99
$wsTasas.Range("O12").Value = 212.522
